$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-79 down to 8-80.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new data record.
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(7, 3).Value = 'Maule'
$ws.Cells.Item(7, 4).Value = 44761
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 100112013
$ws.Cells.Item(7, 7).Value = 'Alcachofa'
$ws.Cells.Item(7, 8).Value = 'Madrigal'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 200
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 16000
$ws.Cells.Item(7, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(7, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 16).Value = 320
$ws.Cells.Item(7, 17).Value = 50
$ws.Cells.Item(7, 18).Value = 'Hortaliza'
